$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix shared string text: remove line break in climate3 label (A10)
$ws.Range("A10").Value = "climate3: Ban the sale of new combustion-engine cars by 2030"

# Update column B values (rows 2-16) with corrected precision figures
$ws.Range("B2").Value = 0.637372904029895
$ws.Range("B3").Value = 0.832071176976898
$ws.Range("B4").Value = 0.806649693500696
$ws.Range("B5").Value = 0.876277717753633
$ws.Range("B6").Value = 0.568918325638433
$ws.Range("B7").Value = 0.583722012843081
$ws.Range("B8").Value = 0.700355311843266
$ws.Range("B9").Value = 0.794665836663391
$ws.Range("B10").Value = 0.624076284194441
$ws.Range("B11").Value = 0.672450987190444
$ws.Range("B12").Value = 0.801839174799188
$ws.Range("B13").Value = 0.711952631580924
$ws.Range("B14").Value = 0.798902933525609
$ws.Range("B15").Value = 0.71470528720425
$ws.Range("B16").Value = 0.582002503082714
